$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header strings (issue number and report week dates) ---
$ws.Range("A8").Value = "Volume 33   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/26/2026  Through  2/1/2026"

# --- Simple numeric value updates ---
$ws.Range("G14").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 20
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -37.5
$ws.Range("F16").Value = 22
$ws.Range("H16").Value = -38.888888888888
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 43
$ws.Range("K16").Value = -41.860465116279
$ws.Range("L16").Value = -30.555555555555
$ws.Range("M16").Value = -16.666666666666
$ws.Range("N16").Value = -86.41304347826
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -30.76923076923
$ws.Range("F17").Value = 52
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 23.809523809523
$ws.Range("I17").Value = 60
$ws.Range("J17").Value = 48
$ws.Range("K17").Value = 25
$ws.Range("L17").Value = 22.448979591836
$ws.Range("M17").Value = 62.162162162162
$ws.Range("N17").Value = -17.808219178082
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -62.962962962963
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 33
$ws.Range("K18").Value = -69.696969696969
$ws.Range("L18").Value = -62.962962962963
$ws.Range("M18").Value = -77.272727272727
$ws.Range("N18").Value = -92.424242424242
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 60
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -22.077922077922
$ws.Range("I19").Value = 63
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = -21.25
$ws.Range("L19").Value = -27.586206896551
$ws.Range("M19").Value = 70.27027027027
$ws.Range("N19").Value = 23.529411764705
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 40
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 36
$ws.Range("J20").Value = 37
$ws.Range("K20").Value = -2.702702702702
$ws.Range("L20").Value = -16.279069767441
$ws.Range("M20").Value = 157.142857142857
$ws.Range("N20").Value = -82.089552238806
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 184
$ws.Range("G21").Value = 221
$ws.Range("H21").Value = -16.742081447963
$ws.Range("I21").Value = 200
$ws.Range("J21").Value = 247
$ws.Range("K21").Value = -19.028340080971
$ws.Range("L21").Value = -18.367346938775
$ws.Range("M21").Value = 21.212121212121
$ws.Range("N21").Value = -69.418960244648
$ws.Range("H22").Value = 50
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = 50
$ws.Range("M22").Value = 50
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = 7.142857142857
$ws.Range("I23").Value = 17
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = 13.333333333333
$ws.Range("L23").Value = -29.166666666666
$ws.Range("M23").Value = 21.428571428571
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -10
$ws.Range("F24").Value = 145
$ws.Range("G24").Value = 165
$ws.Range("H24").Value = -12.121212121212
$ws.Range("I24").Value = 166
$ws.Range("J24").Value = 177
$ws.Range("K24").Value = -6.214689265536
$ws.Range("L24").Value = 0.60606060606
$ws.Range("M24").Value = 38.333333333333
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -30.76923076923
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = -31.372549019607
$ws.Range("I25").Value = 36
$ws.Range("J25").Value = 54
$ws.Range("K25").Value = -33.333333333333
$ws.Range("L25").Value = -38.983050847457
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 240
$ws.Range("F26").Value = 85
$ws.Range("G26").Value = 57
$ws.Range("H26").Value = 49.122807017543
$ws.Range("I26").Value = 96
$ws.Range("J26").Value = 71
$ws.Range("K26").Value = 35.211267605633
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -19.327731092437
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 75
$ws.Range("L27").Value = 40
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("L28").Value = -63.636363636363
$ws.Range("G29").Value = 2
$ws.Range("G30").Value = 2

# --- Cells changing from text placeholder to a real number (copy numeric style, then set value) ---
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("G22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 3
$ws.Range("G22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value = 3
$ws.Range("J22").Copy()
$ws.Range("I22").PasteSpecial(-4122)
$ws.Range("I22").Value = 3
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

# --- Cells changing from a real number to a text placeholder ("0" or "***.*") ---
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
